$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 3 through 15 (the extra BGP-config parameter rows are no
# longer needed - only the "peer"/"type" header rows remain).
$ws.Range("A3:C15").EntireRow.Delete()

# Rename the remaining "peer" row label to "peer-as".
$ws.Range("B2").Value = "peer-as"
